# CSSECDV - Case Study 1 Project Documentation.xlsx
#
# The author widened column D (the "Possible Vulnerabilities" column) on
# Sheet1, which caused a handful of wrapped-text rows to grow taller to
# keep fitting their content, and also scrolled/selected a different cell
# in the frozen-pane view of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column D ("Possible Vulnerabilities") from its previous best-fit
# width to a custom width (~50.125 characters).
$ws.Columns.Item(4).ColumnWidth = 49.33

# The rows whose wrapped text reflows at the new column width need their
# row heights updated (and marked as explicit/custom heights) so the text
# keeps fitting.
$ws.Rows.Item(3).RowHeight = 108.75
$ws.Rows.Item(6).RowHeight = 113.25
$ws.Rows.Item(7).RowHeight = 191.25

# Update the active selection/view to cell C3 (within the frozen-pane
# area), matching the saved view state of the workbook.
$ws.Range("C3").Select()
